$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-case identity cells (B2/B3) and the new User Story cell (B4)
$ws.Range("B2").Value = "verify 40V Load On Addition Deletion of Zetfast loop and its devices "
$ws.Range("B3").Value = "verify40VLoadOnZetfastLoopAddDelete"
$ws.Range("B4").Style = "Normal"
$ws.Range("B4").Value = "NGC-601/T1460 OR TC-211"

# Row 6 header: "Gallery Type" -> "Loops"
$ws.Range("H6").Value = "Loops"

# Expected value text changes "0.00" -> "0.000" (stored as text, not a number)
$ws.Range("J7").Value = "'0.000"
$ws.Range("J8").Value = "'0.000"
$ws.Range("J9").Value = "'0.000"

# Clear the CPU Type value in row 8 (becomes blank/quote-prefixed text)
$ws.Range("C8").Value = "'"

# Move the active selection to B4, matching the saved file's cursor position
[void]$ws.Range("B4").Select()
